# Update DM integration test fixture
#
# 1. Give the header row (row 1) of each sheet its own bold cell style
#    (mirrors the 4 near-identical bold "Arial 11" styles that show up
#    in the target workbook's style table -- one per sheet).
# 2. Re-apply the (now slightly wider, bold-text) column widths.
# 3. Update a handful of ID (UUID) values on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: CodeSchemes ---
$ws1 = $wb.Worksheets.Item("CodeSchemes")
$style1 = $wb.Styles.Add("CodeSchemesHeader")
$style1.Font.Bold = $true
$ws1.Range("A1:N1").Style = "CodeSchemesHeader"

$ws1.Columns.Item(1).ColumnWidth = 28.985714285714288
$ws1.Columns.Item(2).ColumnWidth = 17.435714285714283
$ws1.Columns.Item(3).ColumnWidth = 25.685714285714283
$ws1.Columns.Item(4).ColumnWidth = 22.385714285714283
$ws1.Columns.Item(5).ColumnWidth = 14.135714285714286
$ws1.Columns.Item(6).ColumnWidth = 19.085714285714282
$ws1.Columns.Item(7).ColumnWidth = 22.385714285714286
$ws1.Columns.Item(8).ColumnWidth = 19.085714285714282
$ws1.Columns.Item(9).ColumnWidth = 20.735714285714284
$ws1.Columns.Item(10).ColumnWidth = 24.035714285714285
$ws1.Columns.Item(11).ColumnWidth = 19.085714285714282
$ws1.Columns.Item(12).ColumnWidth = 15.785714285714286
$ws1.Columns.Item(13).ColumnWidth = 20.735714285714284
$ws1.Columns.Item(14).ColumnWidth = 27.335714285714282

$ws1.Range("A2").Value = "5ff6f24a-1497-4f88-a061-ba7745b976f3"

# --- Sheet: Codes ---
$ws2 = $wb.Worksheets.Item("Codes")
$style2 = $wb.Styles.Add("CodesHeader")
$style2.Font.Bold = $true
$ws2.Range("A1:J1").Style = "CodesHeader"

$ws2.Columns.Item(1).ColumnWidth = 34.48571428571429
$ws2.Columns.Item(2).ColumnWidth = 17.435714285714283
$ws2.Columns.Item(3).ColumnWidth = 15.785714285714286
$ws2.Columns.Item(4).ColumnWidth = 14.135714285714286
$ws2.Columns.Item(5).ColumnWidth = 15.785714285714286
$ws2.Columns.Item(6).ColumnWidth = 19.085714285714282
$ws2.Columns.Item(7).ColumnWidth = 20.735714285714284
$ws2.Columns.Item(8).ColumnWidth = 24.035714285714285
$ws2.Columns.Item(9).ColumnWidth = 19.085714285714282
$ws2.Columns.Item(10).ColumnWidth = 15.785714285714286

$ws2.Range("A2").Value = "d9f6f365-03b9-43f5-8370-268380353e6e"
$ws2.Range("A3").Value = "501cf0fd-9181-4860-95ac-ce438485b79a"

# --- Sheet: Extensions ---
$ws3 = $wb.Worksheets.Item("Extensions")
$style3 = $wb.Styles.Add("ExtensionsHeader")
$style3.Font.Bold = $true
$ws3.Range("A1:I1").Style = "ExtensionsHeader"

$ws3.Columns.Item(1).ColumnWidth = 28.985714285714288
$ws3.Columns.Item(2).ColumnWidth = 17.435714285714283
$ws3.Columns.Item(3).ColumnWidth = 14.135714285714286
$ws3.Columns.Item(4).ColumnWidth = 24.035714285714285
$ws3.Columns.Item(5).ColumnWidth = 15.785714285714286
$ws3.Columns.Item(6).ColumnWidth = 19.085714285714282
$ws3.Columns.Item(7).ColumnWidth = 19.085714285714282
$ws3.Columns.Item(8).ColumnWidth = 15.785714285714286
$ws3.Columns.Item(9).ColumnWidth = 24.035714285714285

$ws3.Range("A2").Value = "14364d3f-ae3e-4532-85c7-aa4c2d88af74"

# --- Sheet: Members_dpmDimension ---
$ws4 = $wb.Worksheets.Item("Members_dpmDimension")
$style4 = $wb.Styles.Add("MembersDpmDimensionHeader")
$style4.Font.Bold = $true
$ws4.Range("A1:C1").Style = "MembersDpmDimensionHeader"

$ws4.Columns.Item(1).ColumnWidth = 32.285714285714285
$ws4.Columns.Item(2).ColumnWidth = 10.835714285714285
$ws4.Columns.Item(3).ColumnWidth = 30.635714285714283

$ws4.Range("A2").Value = "8d3ba600-dfea-433c-965a-55d15191fc56"
$ws4.Range("A3").Value = "7c41c694-3725-43c4-8577-4eb2169f082a"
